# Weekly update: a new price record (week of 2022-08-17) is added as a new
# row 26 on the "Maracuyá" sheet, pushing all the existing rows (old 26..59)
# down by one (new 27..60). The sheet's used range grows from T59 to T60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 26; Excel shifts rows 26-59 down
# to 27-60 and copies row 25's formatting into the new row (matches the
# date-format style "s=2" already used by the rest of column D).
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = 10
$ws.Range("B26").Value = 'Vega Modelo de Temuco'
$ws.Range("C26").Value = 'La Araucanía'
$ws.Range("D26").Value = 44790
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = 'Fruta'
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = 'Tropicales y subtropicales'
$ws.Range("I26").Value = 100108003
$ws.Range("J26").Value = 'Maracuyá'
$ws.Range("K26").Value = 'Sin especificar'
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 30
$ws.Range("N26").Value = 36000
$ws.Range("O26").Value = 36000
$ws.Range("P26").Value = 36000
$ws.Range("Q26").Value = '$/caja 18 kilos'
$ws.Range("R26").Value = 'Región de Arica y Parinacota'
$ws.Range("S26").Value = 2000
$ws.Range("T26").Value = 18
